$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 279.34616
$ws.Range("I19").Value = 445.1
$ws.Range("J19").Value = 175.75
$ws.Range("K19").Value = 445.1
$ws.Range("L19").Value = 175.75
$ws.Range("M19").Value = -270.1
$ws.Range("N19").Value = -525.75

$ws.Range("H32").Value = 1890.5416
$ws.Range("I32").Value = 2652.818
$ws.Range("J32").Value = 1245.5385
$ws.Range("K32").Value = 2652.818
$ws.Range("L32").Value = 1245.5385
$ws.Range("M32").Value = -2326.818
$ws.Range("N32").Value = -1897.5385

$ws.Range("H33").Value = 27972.527
$ws.Range("I33").Value = 38657.117
$ws.Range("J33").Value = 192.6
$ws.Range("K33").Value = 38657.117
$ws.Range("L33").Value = 192.6
$ws.Range("M33").Value = -38428.117
$ws.Range("N33").Value = -650.6

$ws.Range("H41").Value = 534.2941
$ws.Range("I41").Value = 250
$ws.Range("J41").Value = 621.7692
$ws.Range("K41").Value = 250
$ws.Range("L41").Value = 621.7692
$ws.Range("M41").Value = 190
$ws.Range("N41").Value = -1501.7692

$ws.Range("H62").Value = 2620.3
$ws.Range("I62").Value = 2457.5715
$ws.Range("K62").Value = 2457.5715
$ws.Range("M62").Value = -1833.5715

$ws.Range("H65").Value = 2620.3
$ws.Range("I65").Value = 2457.5715
$ws.Range("K65").Value = 12287.8575
$ws.Range("M65").Value = -9167.8575

$ws.Range("H98").Value = 1439.9375
$ws.Range("I98").Value = 935.93335
$ws.Range("J98").Value = 9000
$ws.Range("K98").Value = 935.93335
$ws.Range("L98").Value = 9000
$ws.Range("M98").Value = 562.06665
$ws.Range("N98").Value = -11996

$ws.Range("H111").Value = 564
$ws.Range("I111").Value = 472.25
$ws.Range("J111").Value = 747.5
$ws.Range("K111").Value = 1416.75
$ws.Range("L111").Value = 2242.5
$ws.Range("M111").Value = 1650.25
$ws.Range("N111").Value = -8376.5

$ws.Range("H113").Value = 2298.4
$ws.Range("I113").Value = 2020.9286
$ws.Range("J113").Value = 3116.2104
$ws.Range("K113").Value = 2020.9286
$ws.Range("L113").Value = 3116.2104
$ws.Range("M113").Value = 1233.0714
$ws.Range("N113").Value = -9624.2104

$ws.Range("H122").Value = 1439.9375
$ws.Range("I122").Value = 935.93335
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 2807.80005
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -357.8000499999998
$ws.Range("N122").Value = -31900

$ws.Range("H125").Value = 10185.818
$ws.Range("I125").Value = 14920.571
$ws.Range("K125").Value = 134285.139
$ws.Range("M125").Value = -131825.139

$ws.Range("H132").Value = 5440815.5
$ws.Range("I132").Value = 8341784.5
$ws.Range("J132").Value = 1499.5
$ws.Range("K132").Value = 25025353.5
$ws.Range("L132").Value = 4498.5
$ws.Range("M132").Value = -25022823.5
$ws.Range("N132").Value = -9558.5

$ws.Range("H134").Value = 62411.25
$ws.Range("J134").Value = 62411.25
$ws.Range("L134").Value = 62411.25
$ws.Range("N134").Value = -72551.25

$ws.Range("H137").Value = 5171.9243
$ws.Range("I137").Value = 7588.2085
$ws.Range("K137").Value = 22764.6255
$ws.Range("M137").Value = -20214.6255

$ws.Range("H138").Value = 10641796
$ws.Range("I138").Value = 1803.5
$ws.Range("J138").Value = 14289794
$ws.Range("K138").Value = 5410.5
$ws.Range("L138").Value = 42869382
$ws.Range("M138").Value = -270.5
$ws.Range("N138").Value = -42879662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 802.3077
$ws.Range("I2").Value = 681.7
$ws.Range("J2").Value = 1204.3334
$ws.Range("K2").Value = 681.7
$ws.Range("L2").Value = 1204.3334
$ws.Range("M2").Value = -568.7
$ws.Range("N2").Value = -1430.3334

$ws.Range("H45").Value = 1761.375
$ws.Range("I45").Value = 1778.8
$ws.Range("K45").Value = 1778.8
$ws.Range("M45").Value = -1401.8

$ws.Range("H110").Value = 707.75
$ws.Range("I110").Value = 820.1818
$ws.Range("K110").Value = 820.1818
$ws.Range("M110").Value = 1224.8182

$ws.Range("H116").Value = 802.3077
$ws.Range("I116").Value = 681.7
$ws.Range("J116").Value = 1204.3334
$ws.Range("K116").Value = 681.7
$ws.Range("L116").Value = 1204.3334
$ws.Range("M116").Value = 1612.3
$ws.Range("N116").Value = -5792.3334

$ws.Range("H132").Value = 149784.5
$ws.Range("I132").Value = 3913.7256
$ws.Range("J132").Value = 563085
$ws.Range("K132").Value = 11741.1768
$ws.Range("L132").Value = 1689255
$ws.Range("M132").Value = -9211.176800000001
$ws.Range("N132").Value = -1694315

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 802.3077
$ws.Range("I3").Value = 681.7
$ws.Range("J3").Value = 1204.3334
$ws.Range("K3").Value = 681.7
$ws.Range("L3").Value = 1204.3334
$ws.Range("M3").Value = -567.7
$ws.Range("N3").Value = -1432.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15091.384
$ws.Range("I31").Value = 32581.39
$ws.Range("K31").Value = 32581.39
$ws.Range("M31").Value = -32286.39

$ws.Range("H34").Value = 15091.384
$ws.Range("I34").Value = 32581.39
$ws.Range("K34").Value = 32581.39
$ws.Range("M34").Value = -32379.39

$ws.Range("H58").Value = 59338.223
$ws.Range("I58").Value = 5479.4546
$ws.Range("K58").Value = 5479.4546
$ws.Range("M58").Value = -5276.4546

$ws.Range("H132").Value = 3157.1133
$ws.Range("I132").Value = 3444
$ws.Range("K132").Value = 10332
$ws.Range("M132").Value = -7802

$ws.Range("H134").Value = 8951.138999999999
$ws.Range("I134").Value = 6284.696
$ws.Range("K134").Value = 18854.088
$ws.Range("M134").Value = -16319.088

$ws.Range("H136").Value = 59338.223
$ws.Range("I136").Value = 5479.4546
$ws.Range("K136").Value = 16438.3638
$ws.Range("M136").Value = -13888.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1008.7895
$ws.Range("I113").Value = 1363.5714
$ws.Range("J113").Value = 570.5294
$ws.Range("K113").Value = 4090.7142
$ws.Range("L113").Value = 1711.5882
$ws.Range("M113").Value = -1920.7142
$ws.Range("N113").Value = -6051.5882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2375
$ws.Range("I80").Value = 2250
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2250
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1252
$ws.Range("N80").Value = -4496

$ws.Range("H83").Value = 2375
$ws.Range("I83").Value = 2250
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 11250
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -6258
$ws.Range("N83").Value = -22484

$ws.Range("H102").Value = 1200.0667
$ws.Range("I102").Value = 1200.0667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1200.0667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 421.9332999999999
$ws.Range("N102").ClearContents()

$ws.Range("H126").Value = 1917.8572
$ws.Range("J126").Value = 2122.6
$ws.Range("L126").Value = 6367.799999999999
$ws.Range("N126").Value = -11307.8

$ws.Range("H132").Value = 7682.0605
$ws.Range("I132").Value = 6177
$ws.Range("J132").Value = 9997.538
$ws.Range("K132").Value = 18531
$ws.Range("L132").Value = 29992.614
$ws.Range("M132").Value = -16001
$ws.Range("N132").Value = -35052.614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2099.238
$ws.Range("I7").Value = 1907.6666
$ws.Range("J7").Value = 2354.6667
$ws.Range("K7").Value = 1907.6666
$ws.Range("L7").Value = 2354.6667
$ws.Range("M7").Value = -1795.6666
$ws.Range("N7").Value = -2578.6667

$ws.Range("H46").Value = 833.625
$ws.Range("I46").Value = 708.5833
$ws.Range("J46").Value = 958.6667
$ws.Range("K46").Value = 708.5833
$ws.Range("L46").Value = 958.6667
$ws.Range("M46").Value = -520.5833
$ws.Range("N46").Value = -1334.6667

$ws.Range("H126").Value = 2099.238
$ws.Range("I126").Value = 1907.6666
$ws.Range("J126").Value = 2354.6667
$ws.Range("K126").Value = 5722.9998
$ws.Range("L126").Value = 7064.000100000001
$ws.Range("M126").Value = -3252.9998
$ws.Range("N126").Value = -12004.0001

$ws.Range("H132").Value = 22280.152
$ws.Range("I132").Value = 29867.076
$ws.Range("J132").Value = 7485.65
$ws.Range("K132").Value = 89601.228
$ws.Range("L132").Value = 22456.95
$ws.Range("M132").Value = -87071.228
$ws.Range("N132").Value = -27516.95

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 299.41666
$ws.Range("I107").Value = 143.66667
$ws.Range("J107").Value = 766.6667
$ws.Range("K107").Value = 431.00001
$ws.Range("L107").Value = 2300.0001
$ws.Range("M107").Value = 1488.99999
$ws.Range("N107").Value = -6140.0001

$ws.Range("H126").Value = 501156.2
$ws.Range("I126").Value = 909988.0600000001
$ws.Range("J126").Value = 1472.7778
$ws.Range("K126").Value = 2729964.18
$ws.Range("L126").Value = 4418.3334
$ws.Range("M126").Value = -2727494.18
$ws.Range("N126").Value = -9358.3334

$ws.Range("H132").Value = 7035
$ws.Range("I132").Value = 8120.6333
$ws.Range("J132").Value = 4863.7334
$ws.Range("K132").Value = 24361.8999
$ws.Range("L132").Value = 14591.2002
$ws.Range("M132").Value = -21831.8999
$ws.Range("N132").Value = -19651.2002

$ws.Range("H136").Value = 21742596
$ws.Range("I136").Value = 32262298
$ws.Range("J136").Value = 1878.1333
$ws.Range("K136").Value = 96786894
$ws.Range("L136").Value = 5634.3999
$ws.Range("M136").Value = -96784344
$ws.Range("N136").Value = -10734.3999
